# Lab01 Review Report update
# - Fills in the "Requirements Phase Defects" sheet with the reviewer's
#   info and three logged defects (R01-R03), and stamps the reviewer
#   name / review date onto all three review-form sheets, plus records
#   a 0.5h review effort on each sheet. Also removes the old
#   "do not print this form" placeholder text from B1 on every sheet
#   and fixes the casing of the Coding sheet's document-title value.

$wb = $excel.ActiveWorkbook

$reviewer   = "Tritean Tudor-Adrian"
$reviewDate = "10.03.2019"

$sheetNames = @(
    "Requirements Phase Defects",
    "Architect. Design Phase Defects",
    "Coding Phase Defects"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Remove the "do not print this form" watermark text in B1.
    $null = $ws.Range("B1").ClearContents()

    # Stamp reviewer name (Author Name value cell, D5 on every sheet).
    $ws.Range("D5").Value = $reviewer

    # Stamp the review date as literal text (not an Excel date serial).
    $ws.Range("D7").NumberFormat = "@"
    $ws.Range("D7").Value = $reviewDate

    # Record effort spent reviewing (last row, column E) as 0.5 hours.
    $dim = $ws.UsedRange
    $lastRow = $dim.Row + $dim.Rows.Count - 1
    $ws.Cells.Item($lastRow, 5).Value = 0.5
}

# Coding sheet: fix document title casing ("Coding Document" -> "Coding document").
$wsCoding = $wb.Worksheets.Item("Coding Phase Defects")
$wsCoding.Range("D4").Value = "Coding document"

# Requirements sheet: log the three defects found while reviewing.
$wsReq = $wb.Worksheets.Item("Requirements Phase Defects")
$wsReq.Range("C10").Value = "R01"
$wsReq.Range("E10").Value = "The type of the report is not defined(pdf,etc.)"
$wsReq.Range("C11").Value = "R02"
$wsReq.Range("E11").Value = "Authentication is not even mentioned"
$wsReq.Range("C12").Value = "R03"
$wsReq.Range("E12").Value = "Properties of a user are defined while describing a functionality"
$wsReq.Rows.Item(12).RowHeight = 30

# Selections left behind in each sheet + make the Requirements sheet active/selected.
$wsArch = $wb.Worksheets.Item("Architect. Design Phase Defects")
$null = $wsArch.Range("E28").Select()

$null = $wsCoding.Range("E32").Select()

$wsReq.Activate()
$null = $wsReq.Range("E13").Select()
